$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New working set of sequences (image path, German verb, category) replacing
# the previous rows 2-33 (data indices 0-31). Column A (index) and header
# row are left untouched.
$data = @(
    @(2, 62, "flower/flower006.jpg", "lehnen", "flower"),
    @(3, 10, "dog/dog024.jpg", "töten", "dog"),
    @(4, 94, "dog/dog009.jpg", "regnen", "dog"),
    @(5, 66, "dog/dog014.jpg", "sondern", "dog"),
    @(6, 13, "dog/dog002.jpg", "drohen", "dog"),
    @(7, 90, "dog/dog000.jpg", "schmecken", "dog"),
    @(8, 103, "flower/flower027.jpg", "biegen", "flower"),
    @(9, 48, "flower/flower025.jpg", "formen", "flower"),
    @(10, 16, "flower/flower013.jpg", "wiegen", "flower"),
    @(11, 25, "dog/dog028.jpg", "langen", "dog"),
    @(12, 78, "flower/flower022.jpg", "ehren", "flower"),
    @(13, 99, "flower/flower021.jpg", "posten", "flower"),
    @(14, 113, "flower/flower017.jpg", "starten", "flower"),
    @(15, 114, "dog/dog021.jpg", "fühlen", "dog"),
    @(16, 7, "dog/dog027.jpg", "spielen", "dog"),
    @(17, 105, "flower/flower010.jpg", "klappen", "flower"),
    @(18, 45, "dog/dog011.jpg", "sieben", "dog"),
    @(19, 0, "flower/flower019.jpg", "bitten", "flower"),
    @(20, 55, "flower/flower003.jpg", "strahlen", "flower"),
    @(21, 77, "flower/flower014.jpg", "kehren", "flower"),
    @(22, 84, "dog/dog029.jpg", "mieten", "dog"),
    @(23, 95, "flower/flower026.jpg", "jubeln", "flower"),
    @(24, 57, "flower/flower002.jpg", "krachen", "flower"),
    @(25, 70, "dog/dog012.jpg", "hoffen", "dog"),
    @(26, 40, "flower/flower000.jpg", "tauschen", "flower"),
    @(27, 92, "dog/dog007.jpg", "füttern", "dog"),
    @(28, 83, "dog/dog026.jpg", "währen", "dog"),
    @(29, 4, "flower/flower008.jpg", "pflegen", "flower"),
    @(30, 21, "dog/dog025.jpg", "füllen", "dog"),
    @(31, 93, "flower/flower028.jpg", "schätzen", "flower"),
    @(32, 28, "dog/dog001.jpg", "runden", "dog"),
    @(33, 36, "dog/dog003.jpg", "haken", "dog")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
